$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "FilesTab" query in cell B4: drop the `File Type` and `Breed`
# columns from the RETURN clause (Bento object repository revisited).
$newQuery = "`nMATCH (f:file)-->(parent)`nWITH DISTINCT f, parent`nMATCH (f)-[*]->(c:case)<--(demo:demographic)`n MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)`nWHERE demo.sex IN ['Unknown']`nWITH DISTINCT f, parent, c, demo, diag, s`nRETURN coalesce(f.file_name, '') AS ``File Name``, `n        coalesce(labels(parent)[0], '') AS ``Association``,`n        coalesce(f.file_description, '') AS ``Description``,`n        coalesce(f.file_format, '') AS ``Format``,`n        coalesce(f.file_size, '') AS ``Size``,`n        coalesce(c.case_id, '') AS ``Case ID``, `n        coalesce(diag.disease_term,'') AS Diagnosis , `n        coalesce(s.clinical_study_designation,'') AS ``Study Code``"

$ws.Range("B4").Value = $newQuery

# Match the author's scroll/selection change: the sheet view now shows row 4
# at the top with B4 selected instead of B2.
$ws.Range("B4").Select()
$excel.ActiveWindow.ScrollRow = 4
